$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 changes
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8

# Row 11 changes
$ws.Range("G11").Value = 4.2
$ws.Range("I11").Value = 1.9
$ws.Range("L11").Value = 2.63
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("W11").Value = 9
$ws.Range("AE11").Value = 19
$ws.Range("AI11").Value = 8
$ws.Range("AJ11").Value = 9.5
$ws.Range("BB11").Value = 251
